$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "MIT 6.800/6.843 Robotics Manipulation :: Motion planning(1)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/187"

$ws.Range("D32").Value = "Optimizers 개념 정리"
$ws.Range("E32").Value = "https://dodonam.tistory.com/356"

$ws.Range("D51").Value = "[python+pandas] 판다스 데이터 프레임에서 컬럼의 고유값을 알고 싶으면, unique 메소드"
$ws.Range("E51").Value = "https://bskyvision.com/1267"
